$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update placeholder cell values ---
# The "Total X:" amount cells (Q3:Q12) used to hold literal "[TotalXXX]"
# placeholder tokens; the refreshed template now starts them at 0 so the
# new named ranges below can be wired up to real formulas/macros.
$ws.Range("Q3").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("Q12").Value = 0

# Row 16 was a sample/placeholder data row ([Size], [SO], [DesignNum],
# [Batch], [Quantity], [SchNumber]). Remove all of that sample content;
# A16 now becomes the (empty) anchor cell for the named range below, and
# P16 is simply cleared while keeping its existing formatting.
$ws.Range("A16").Clear()
$ws.Range("C16").Clear()
$ws.Range("E16").Clear()
$ws.Range("G16").Clear()
$ws.Range("I16").Clear()
$ws.Range("P16").ClearContents()

# Put the selection where the template now wants it.
$ws.Range("A16").Select()

# --- Replace the old global DataRange name with sheet-scoped named cells ---
$wb.Names.Item("DataRange").Delete()

$ws.Names.Add('FirstCellOfPacklistLineData', '=PackingList!$A$16')
$ws.Names.Add('GrandTotal', '=PackingList!$Q$12')
$ws.Names.Add('Total10WAV', '=PackingList!$Q$6')
$ws.Names.Add('Total16BEER', '=PackingList!$Q$9')
$ws.Names.Add('Total16DWT', '=PackingList!$Q$7')
$ws.Names.Add('Total16MUG', '=PackingList!$Q$8')
$ws.Names.Add('Total24DWT', '=PackingList!$Q$10')
$ws.Names.Add('Total24WB', '=PackingList!$Q$11')
$ws.Names.Add('Total6SIP', '=PackingList!$Q$3')
$ws.Names.Add('Total9SWG', '=PackingList!$Q$4')
$ws.Names.Add('Total9WINE', '=PackingList!$Q$5')
